# Apply the "registros.xlsx" -> Clientes sheet fix:
#  - resize columns A:C
#  - update the sample record in row 2 (new cedula, new name, new telefono)
#  - clear the "Activo" flag cell (column D) for that row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Excel's ColumnWidth property is expressed in characters of the Normal
# style's font, which carries an implicit ~0.8333 padding offset versus the
# raw width stored in the OOXML <col> element. Subtract that offset so the
# saved width matches the target values exactly (11, 17, 12).
$padding = 0.8333333333333333
$ws.Columns.Item(1).ColumnWidth = 11 - $padding
$ws.Columns.Item(2).ColumnWidth = 17 - $padding
$ws.Columns.Item(3).ColumnWidth = 12 - $padding

# Update row 2 values
$ws.Range("A2").Value = 901234359
$ws.Range("B2").Value = "Pedro la piedra"
$ws.Range("C2").Value = 3453682345

# Remove the "Activo" boolean value entirely for this row
$ws.Range("D2").ClearContents()
